# fix: excel for overtime
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header labels from ATTENDANCE_* to OVERTIME_*
$ws.Range("B1").Value = "OVERTIME_DATE"
$ws.Range("C1").Value = "OVERTIME_IN"
$ws.Range("D1").Value = "OVERTIME_OUT"

# Update the selected/active cell shown when the workbook is opened
$ws.Range("E4").Select()
